$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-8), columns: D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), P (Precio $/Kg), Q (Kg o Unidades)

$data = @{
    2 = @{ D = 44293; I = "Primera"; J = 10; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos empedrada"; P = 1667; Q = 15 }
    3 = @{ D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    4 = @{ D = 44315; I = "Primera"; J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel"; P = 1000; Q = 15 }
    5 = @{ D = 44280; I = "Primera"; J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    6 = @{ D = 44313; I = "Primera"; J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 }
    7 = @{ D = 44313; I = "Primera"; J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    8 = @{ D = 44285; I = "Primera"; J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("I$row").Value = $vals.I
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
}
